$d = $word.ActiveDocument

# Find the "LOB1006: Cálculo IV (Requisito)" paragraph; its end (including the
# paragraph mark) is where the block to remove begins.
$startRange = $d.Content
[void]$startRange.Find.Execute("LOB1006: Cálculo IV (Requisito)", $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
[void]$startRange.Expand(4)  # wdParagraph - include the paragraph mark

# Find the "© 2020 . Contact: ..." paragraph; its end (including the paragraph
# mark) is where the block to remove ends.
$endRange = $d.Content
[void]$endRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
[void]$endRange.Expand(4)  # wdParagraph - include the paragraph mark

# Delete everything from right after "LOB1006: Cálculo IV (Requisito)" through
# the end of the "© 2020 ..." paragraph. This removes the blank paragraph
# between them, the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph,
# and the "© 2020 ..." paragraph itself, while leaving one blank paragraph
# (the one that used to follow the copyright line) in place before the
# trailing page-break paragraph.
$deleteRange = $d.Range($startRange.End, $endRange.End)
$deleteRange.Delete()
